# Update test data for removed variables and rely more on moorings &
# installation examples.
#
# On the "On-Site" sheet, row 3 (previously blank) is populated with a new
# "Mooring Lines" sub-system entry, mirroring the existing "Foundations" row.
# A new shared string "Mooring Lines" is introduced as a side effect of
# setting the cell's text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("On-Site")

# Column A on row 3 should pick up the same number format/style as the other
# data cells in that row (B3:E3 use style index 1, same as row 2's data
# cells), rather than the style used by the other blank column-A cells.
# Copy the formatting from B3 (a plain data cell) onto A3 before writing the
# label, then fill in the rest of the row's values.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A3").Value = "Mooring Lines"
$ws.Range("B3").Value = 5000
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 3

# Make "On-Site" the active sheet with A3 selected (it was previously the
# "Replacement" sheet that held the active/selected state).
$ws.Activate() | Out-Null
$ws.Range("A3").Select() | Out-Null
